$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A14: change from text to numeric value
$ws.Cells.Item(14, 1).Value = 79174445

# Add new row 15
$ws.Cells.Item(15, 1).Value = 79174445
$ws.Cells.Item(15, 2).Value = ""
$ws.Cells.Item(15, 3).Value = "Cash"
$ws.Cells.Item(15, 4).Value = "2025-08-18T08:46:52"
$ws.Cells.Item(15, 5).Value = 30
$ws.Cells.Item(15, 6).Value = ""
$ws.Cells.Item(15, 7).Value = 10
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 20
